$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 375.8889
$ws.Cells.Item(6, 9).Value = 48
$ws.Cells.Item(6, 11).Value = 144
$ws.Cells.Item(6, 13).Value = -32
$ws.Cells.Item(8, 8).Value = 1089.3334
$ws.Cells.Item(8, 9).Value = 160.4
$ws.Cells.Item(8, 10).Value = 2250.5
$ws.Cells.Item(8, 11).Value = 481.2
$ws.Cells.Item(8, 12).Value = 6751.5
$ws.Cells.Item(8, 13).Value = -342.2
$ws.Cells.Item(8, 14).Value = -7029.5
$ws.Cells.Item(43, 8).Value = 500
$ws.Cells.Item(43, 10).Value = 500
$ws.Cells.Item(43, 12).Value = 500
$ws.Cells.Item(43, 14).Value = -638
$ws.Cells.Item(99, 8).Value = 1785.3334
$ws.Cells.Item(99, 9).Value = 138
$ws.Cells.Item(99, 10).Value = 5080
$ws.Cells.Item(99, 11).Value = 414
$ws.Cells.Item(99, 12).Value = 15240
$ws.Cells.Item(99, 13).Value = 1084
$ws.Cells.Item(99, 14).Value = -18236
$ws.Cells.Item(135, 8).Value = 1314.5
$ws.Cells.Item(135, 9).Value = 1195.3
$ws.Cells.Item(135, 10).Value = 1910.5
$ws.Cells.Item(135, 11).Value = 10757.7
$ws.Cells.Item(135, 12).Value = 17194.5
$ws.Cells.Item(135, 13).Value = -8222.699999999999
$ws.Cells.Item(135, 14).Value = -22264.5
$ws.Cells.Item(137, 8).Value = 1688.4839
$ws.Cells.Item(137, 10).Value = 2690.6155
$ws.Cells.Item(137, 12).Value = 8071.8465
$ws.Cells.Item(137, 14).Value = -13171.8465
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 301
$ws.Cells.Item(5, 9).Value = 301
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 301
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = -189
$ws.Cells.Item(5, 14).ClearContents()
$ws.Cells.Item(32, 8).Value = 9093630
$ws.Cells.Item(32, 9).Value = 2992.6
$ws.Cells.Item(32, 11).Value = 2992.6
$ws.Cells.Item(32, 13).Value = -2705.6
$ws.Cells.Item(41, 8).Value = 3028
$ws.Cells.Item(41, 9).Value = 556
$ws.Cells.Item(41, 10).Value = 5500
$ws.Cells.Item(41, 11).Value = 556
$ws.Cells.Item(41, 12).Value = 5500
$ws.Cells.Item(41, 13).Value = -142
$ws.Cells.Item(41, 14).Value = -6328
$ws.Cells.Item(61, 8).Value = 3010.6428
$ws.Cells.Item(61, 9).Value = 2901.3333
$ws.Cells.Item(61, 10).Value = 3666.5
$ws.Cells.Item(61, 11).Value = 2901.3333
$ws.Cells.Item(61, 12).Value = 3666.5
$ws.Cells.Item(61, 13).Value = -2689.3333
$ws.Cells.Item(61, 14).Value = -4090.5
$ws.Cells.Item(74, 8).Value = 2193.5
$ws.Cells.Item(74, 9).Value = 1791.3
$ws.Cells.Item(74, 11).Value = 1791.3
$ws.Cells.Item(74, 13).Value = -917.3
$ws.Cells.Item(77, 8).Value = 2193.5
$ws.Cells.Item(77, 9).Value = 1791.3
$ws.Cells.Item(77, 11).Value = 8956.5
$ws.Cells.Item(77, 13).Value = -4588.5
$ws.Cells.Item(97, 8).Value = 844.44446
$ws.Cells.Item(97, 9).Value = 844.44446
$ws.Cells.Item(97, 11).Value = 844.44446
$ws.Cells.Item(97, 13).Value = -348.44446
$ws.Cells.Item(132, 8).Value = 1632.625
$ws.Cells.Item(132, 9).Value = 1632.625
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 4897.875
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -2367.875
$ws.Cells.Item(132, 14).ClearContents()
$ws.Cells.Item(136, 8).Value = 3010.6428
$ws.Cells.Item(136, 9).Value = 2901.3333
$ws.Cells.Item(136, 10).Value = 3666.5
$ws.Cells.Item(136, 11).Value = 8703.999899999999
$ws.Cells.Item(136, 12).Value = 10999.5
$ws.Cells.Item(136, 13).Value = -6153.999899999999
$ws.Cells.Item(136, 14).Value = -16099.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 301
$ws.Cells.Item(4, 9).Value = 301
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 301
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = -186
$ws.Cells.Item(4, 14).ClearContents()
$ws.Cells.Item(36, 8).Value = 988.8570999999999
$ws.Cells.Item(36, 9).Value = 988.8570999999999
$ws.Cells.Item(36, 11).Value = 988.8570999999999
$ws.Cells.Item(36, 13).Value = -454.8570999999999
$ws.Cells.Item(86, 8).Value = 6602.0835
$ws.Cells.Item(86, 9).Value = 3445
$ws.Cells.Item(86, 11).Value = 3445
$ws.Cells.Item(86, 13).Value = -2322
$ws.Cells.Item(89, 8).Value = 6602.0835
$ws.Cells.Item(89, 9).Value = 3445
$ws.Cells.Item(89, 11).Value = 17225
$ws.Cells.Item(89, 13).Value = -11609
$ws.Cells.Item(134, 8).Value = 3715.682
$ws.Cells.Item(134, 9).Value = 876.1053000000001
$ws.Cells.Item(134, 11).Value = 2628.3159
$ws.Cells.Item(134, 13).Value = -93.31590000000006
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(56, 8).Value = 0
$ws.Cells.Item(56, 9).Value = 0
$ws.Cells.Item(56, 11).Value = 0
$ws.Cells.Item(56, 13).ClearContents()
$ws.Cells.Item(58, 8).Value = 8577
$ws.Cells.Item(58, 9).Value = 7997.5
$ws.Cells.Item(58, 11).Value = 7997.5
$ws.Cells.Item(58, 13).Value = -7794.5
$ws.Cells.Item(62, 8).Value = 2773.1667
$ws.Cells.Item(62, 9).Value = 2450
$ws.Cells.Item(62, 10).Value = 2934.75
$ws.Cells.Item(62, 11).Value = 2450
$ws.Cells.Item(62, 12).Value = 2934.75
$ws.Cells.Item(62, 13).Value = -1826
$ws.Cells.Item(62, 14).Value = -4182.75
$ws.Cells.Item(65, 8).Value = 2773.1667
$ws.Cells.Item(65, 9).Value = 2450
$ws.Cells.Item(65, 10).Value = 2934.75
$ws.Cells.Item(65, 11).Value = 12250
$ws.Cells.Item(65, 12).Value = 14673.75
$ws.Cells.Item(65, 13).Value = -9130
$ws.Cells.Item(65, 14).Value = -20913.75
$ws.Cells.Item(134, 8).Value = 3905.4167
$ws.Cells.Item(134, 9).Value = 3816.0908
$ws.Cells.Item(134, 11).Value = 11448.2724
$ws.Cells.Item(134, 13).Value = -8913.2724
$ws.Cells.Item(136, 8).Value = 8577
$ws.Cells.Item(136, 9).Value = 7997.5
$ws.Cells.Item(136, 11).Value = 23992.5
$ws.Cells.Item(136, 13).Value = -21442.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 947.125
$ws.Cells.Item(17, 9).Value = 694.2
$ws.Cells.Item(17, 10).Value = 1368.6666
$ws.Cells.Item(17, 11).Value = 2082.6
$ws.Cells.Item(17, 12).Value = 4105.9998
$ws.Cells.Item(17, 13).Value = -1913.6
$ws.Cells.Item(17, 14).Value = -4443.9998
$ws.Cells.Item(34, 8).Value = 746.7778
$ws.Cells.Item(34, 9).Value = 185.4
$ws.Cells.Item(34, 10).Value = 962.6923
$ws.Cells.Item(34, 11).Value = 556.2
$ws.Cells.Item(34, 12).Value = 2888.0769
$ws.Cells.Item(34, 13).Value = -472.2
$ws.Cells.Item(34, 14).Value = -3056.0769
$ws.Cells.Item(40, 8).Value = 178.81818
$ws.Cells.Item(40, 10).Value = 363.8
$ws.Cells.Item(40, 12).Value = 1455.2
$ws.Cells.Item(40, 14).Value = -1593.2
$ws.Cells.Item(49, 8).Value = 2987.5
$ws.Cells.Item(49, 9).Value = 2975
$ws.Cells.Item(49, 11).Value = 8925
$ws.Cells.Item(49, 13).Value = -8769
$ws.Cells.Item(107, 8).Value = 462.6
$ws.Cells.Item(107, 9).Value = 221.83333
$ws.Cells.Item(107, 10).Value = 717.5294
$ws.Cells.Item(107, 11).Value = 665.49999
$ws.Cells.Item(107, 12).Value = 2152.5882
$ws.Cells.Item(107, 13).Value = 1254.50001
$ws.Cells.Item(107, 14).Value = -5992.5882
$ws.Cells.Item(119, 8).Value = 4615.8
$ws.Cells.Item(119, 9).Value = 4615.8
$ws.Cells.Item(119, 11).Value = 13847.4
$ws.Cells.Item(119, 13).Value = -9009.400000000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2145.4443
$ws.Cells.Item(80, 9).Value = 2272.7144
$ws.Cells.Item(80, 10).Value = 1700
$ws.Cells.Item(80, 11).Value = 2272.7144
$ws.Cells.Item(80, 12).Value = 1700
$ws.Cells.Item(80, 13).Value = -1274.7144
$ws.Cells.Item(80, 14).Value = -3696
$ws.Cells.Item(83, 8).Value = 2145.4443
$ws.Cells.Item(83, 9).Value = 2272.7144
$ws.Cells.Item(83, 10).Value = 1700
$ws.Cells.Item(83, 11).Value = 11363.572
$ws.Cells.Item(83, 12).Value = 8500
$ws.Cells.Item(83, 13).Value = -6371.572
$ws.Cells.Item(83, 14).Value = -18484
$ws.Cells.Item(122, 8).Value = 2329.3572
$ws.Cells.Item(122, 9).Value = 1600.7142
$ws.Cells.Item(122, 10).Value = 3058
$ws.Cells.Item(122, 11).Value = 4802.142599999999
$ws.Cells.Item(122, 12).Value = 9174
$ws.Cells.Item(122, 13).Value = -2352.142599999999
$ws.Cells.Item(122, 14).Value = -14074
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2335
$ws.Cells.Item(7, 9).Value = 2335
$ws.Cells.Item(7, 11).Value = 2335
$ws.Cells.Item(7, 13).Value = -2223
$ws.Cells.Item(16, 8).Value = 1703.1
$ws.Cells.Item(16, 9).Value = 1733.1111
$ws.Cells.Item(16, 10).Value = 1433
$ws.Cells.Item(16, 11).Value = 1733.1111
$ws.Cells.Item(16, 12).Value = 1433
$ws.Cells.Item(16, 13).Value = -1563.1111
$ws.Cells.Item(16, 14).Value = -1773
$ws.Cells.Item(61, 8).Value = 5222.5386
$ws.Cells.Item(61, 9).Value = 4049.25
$ws.Cells.Item(61, 11).Value = 4049.25
$ws.Cells.Item(61, 13).Value = -3847.25
$ws.Cells.Item(68, 8).Value = 4061.158
$ws.Cells.Item(68, 9).Value = 1890.3334
$ws.Cells.Item(68, 11).Value = 1890.3334
$ws.Cells.Item(68, 13).Value = -1141.3334
$ws.Cells.Item(71, 8).Value = 4061.158
$ws.Cells.Item(71, 9).Value = 1890.3334
$ws.Cells.Item(71, 11).Value = 9451.666999999999
$ws.Cells.Item(71, 13).Value = -5707.666999999999
$ws.Cells.Item(93, 8).Value = 1000
$ws.Cells.Item(93, 9).Value = 1000
$ws.Cells.Item(93, 11).Value = 1000
$ws.Cells.Item(93, 13).Value = 248
$ws.Cells.Item(100, 8).Value = 8499.857
$ws.Cells.Item(100, 9).Value = 5000
$ws.Cells.Item(100, 11).Value = 5000
$ws.Cells.Item(100, 13).Value = -4459
$ws.Cells.Item(113, 8).Value = 5222.5386
$ws.Cells.Item(113, 9).Value = 4049.25
$ws.Cells.Item(113, 11).Value = 4049.25
$ws.Cells.Item(113, 13).Value = -1879.25
$ws.Cells.Item(126, 8).Value = 2335
$ws.Cells.Item(126, 9).Value = 2335
$ws.Cells.Item(126, 11).Value = 7005
$ws.Cells.Item(126, 13).Value = -4535
$ws.Cells.Item(132, 8).Value = 3377.889
$ws.Cells.Item(132, 9).Value = 3160.2
$ws.Cells.Item(132, 10).Value = 3650
$ws.Cells.Item(132, 11).Value = 9480.599999999999
$ws.Cells.Item(132, 12).Value = 10950
$ws.Cells.Item(132, 13).Value = -6950.599999999999
$ws.Cells.Item(132, 14).Value = -16010
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(58, 8).Value = 3112.4443
$ws.Cells.Item(58, 9).Value = 3112.4443
$ws.Cells.Item(58, 11).Value = 3112.4443
$ws.Cells.Item(58, 13).Value = -2804.4443
$ws.Cells.Item(96, 8).Value = 998.2
$ws.Cells.Item(96, 10).Value = 745
$ws.Cells.Item(96, 12).Value = 745
$ws.Cells.Item(96, 14).Value = -3491
$ws.Cells.Item(113, 8).Value = 846.3333
$ws.Cells.Item(113, 9).Value = 676
$ws.Cells.Item(113, 11).Value = 2028
$ws.Cells.Item(113, 13).Value = 142
$ws.Cells.Item(122, 8).Value = 2745.8635
$ws.Cells.Item(122, 9).Value = 2573.8948
$ws.Cells.Item(122, 11).Value = 7721.6844
$ws.Cells.Item(122, 13).Value = -5271.6844
$ws.Cells.Item(136, 8).Value = 3411.9375
$ws.Cells.Item(136, 9).Value = 2755.6667
$ws.Cells.Item(136, 10).Value = 3805.7
$ws.Cells.Item(136, 11).Value = 8267.000100000001
$ws.Cells.Item(136, 12).Value = 11417.1
$ws.Cells.Item(136, 13).Value = -5717.000100000001
$ws.Cells.Item(136, 14).Value = -16517.1
